$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap full row contents between row 3 (West Virginia Primary) and row 5
# (Clinton haunted by coal country comment), leaving row 2 and row 4 as-is.
# (Matches the reordering seen in the shared-strings table: rows 3 and 5
# traded places while the hyperlink relationships stayed tied to their
# original rows, exactly like the upstream diff shows.)

$a3 = $ws.Range("A3").Value2
$b3 = $ws.Range("B3").Value2
$c3 = $ws.Range("C3").Value2
$d3 = $ws.Range("D3").Value2
$e3 = $ws.Range("E3").Value2

$a5 = $ws.Range("A5").Value2
$b5 = $ws.Range("B5").Value2
$c5 = $ws.Range("C5").Value2
$d5 = $ws.Range("D5").Value2
$e5 = $ws.Range("E5").Value2

$ws.Range("A3").Value = $a5
$ws.Range("B3").Value = $b5
$ws.Range("C3").Value = $c5
$ws.Range("D3").Value = $d5
$ws.Range("E3").Value = $e5

$ws.Range("A5").Value = $a3
$ws.Range("B5").Value = $b3
$ws.Range("C5").Value = $c3
$ws.Range("D5").Value = $d3
$ws.Range("E5").Value = $e3
